function Wrap-Xml($inner) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $inner + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Replace-ParagraphRange($range, $newParaXml) {
    $range.Collapse(0)
    $range.MoveEnd(1, -1)
    $range.InsertXML((Wrap-Xml($newParaXml)))
}

$d = $word.ActiveDocument

# 1) "Come up with an idea for this project " -> add run "- Done"
$p = $d.Paragraphs(6)
$xml = '<w:p w14:paraId="713F379E" w14:textId="7B6D6A08" w:rsidR="005C0182" w:rsidRDefault="005C0182" w:rsidP="005C0182">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Come up with an idea for this project </w:t></w:r>' + `
  '<w:r><w:t>- Done</w:t></w:r>' + `
  '</w:p>'
Replace-ParagraphRange $p.Range $xml

# 2) "Create the read me file and submit via blackboard" -> add run " - Done"
$p = $d.Paragraphs(7)
$xml = '<w:p w14:paraId="0CAE70D9" w14:textId="243ECC32" w:rsidR="005C0182" w:rsidRDefault="005C0182" w:rsidP="005C0182">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Create the read me file and submit via blackboard</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> - Done</w:t></w:r>' + `
  '</w:p>'
Replace-ParagraphRange $p.Range $xml

# 3) "Start on the code for the games" -> add run " - Done"
$p = $d.Paragraphs(8)
$xml = '<w:p w14:paraId="182C6C4A" w14:textId="71D1CF4B" w:rsidR="005C0182" w:rsidRDefault="005C0182" w:rsidP="005C0182">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Start on the code for the games</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> - Done</w:t></w:r>' + `
  '</w:p>'
Replace-ParagraphRange $p.Range $xml

# 4) "Set up a meeting date by 12/10" -> add run " - Done"
$p = $d.Paragraphs(9)
$xml = '<w:p w14:paraId="190A6536" w14:textId="4FCC0973" w:rsidR="005C0182" w:rsidRDefault="005C0182" w:rsidP="005C0182">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Set up a meeting date by 12/10</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> - Done</w:t></w:r>' + `
  '</w:p>'
Replace-ParagraphRange $p.Range $xml

# 5) "Work on more code. Try to be 50% by meeting date" -> add run " - Done",
#    then insert two brand-new list paragraphs right after it (before "Submit completed work").
$p = $d.Paragraphs(10)
$xml = '<w:p w14:paraId="61936B53" w14:textId="544F70B4" w:rsidR="005C0182" w:rsidRDefault="005C0182" w:rsidP="005C0182">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Work on more code. Try to be 50% by meeting date</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> - Done</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Upload code onto Git hub – Done</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
  '<w:r><w:t>Meeting on 12/14 @17:45</w:t></w:r></w:p>'
Replace-ParagraphRange $p.Range $xml
